$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the combined string: a Python-tuple-like repr combining the card's
# name and its remaining attributes.
$name = "Ryusei, the Falling Star"
$combined = "('" + $name + "', ['{5}{R}', 'Legendary Creature " + [char]0x2014 + " Dragon Spirit', 'Flying', 'When " + $name + " dies, it deals 5 damage to each creature without flying.', '5/5'])"

# Set the new combined value into A2
$ws.Range("A2").Value = $combined

# Remove the now-unneeded rows (previously A3:A7) that held the separate fields
$ws.Range("A3:A7").ClearContents()
